# Auto-generated edit script: update leve profit data cells per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1874.5
$ws.Range("J70").Value = 2066.6667
$ws.Range("L70").Value = 6200.000100000001
$ws.Range("N70").Value = -6740.000100000001
$ws.Range("H73").Value = 1874.5
$ws.Range("J73").Value = 2066.6667
$ws.Range("L73").Value = 6200.000100000001
$ws.Range("N73").Value = -8072.000100000001
$ws.Range("H135").Value = 377.35715
$ws.Range("I135").Value = 427
$ws.Range("J135").Value = 195.33333
$ws.Range("K135").Value = 3843
$ws.Range("L135").Value = 1757.99997
$ws.Range("M135").Value = -1308
$ws.Range("N135").Value = -6827.99997
$ws.Range("H138").Value = 4198.7896
$ws.Range("J138").Value = 4948
$ws.Range("L138").Value = 14844
$ws.Range("N138").Value = -25124

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2990
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H32").Value = 3211694.2
$ws.Range("I32").Value = 5004349
$ws.Range("J32").Value = 701977.7
$ws.Range("K32").Value = 5004349
$ws.Range("L32").Value = 701977.7
$ws.Range("M32").Value = -5004062
$ws.Range("N32").Value = -702551.7
$ws.Range("H45").Value = 3149.4546
$ws.Range("I45").Value = 3569.889
$ws.Range("K45").Value = 3569.889
$ws.Range("M45").Value = -3192.889
$ws.Range("H74").Value = 3011
$ws.Range("I74").Value = 3011
$ws.Range("K74").Value = 3011
$ws.Range("M74").Value = -2137
$ws.Range("H77").Value = 3011
$ws.Range("I77").Value = 3011
$ws.Range("K77").Value = 15055
$ws.Range("M77").Value = -10687
$ws.Range("H116").Value = 2990
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 1358.9656
$ws.Range("I122").Value = 1226.2916
$ws.Range("J122").Value = 1995.8
$ws.Range("K122").Value = 3678.8748
$ws.Range("L122").Value = 5987.4
$ws.Range("M122").Value = -1228.8748
$ws.Range("N122").Value = -10887.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2990
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H64").Value = 857
$ws.Range("J64").Value = 789.75
$ws.Range("L64").Value = 789.75
$ws.Range("N64").Value = -1239.75
$ws.Range("H67").Value = 857
$ws.Range("J67").Value = 789.75
$ws.Range("L67").Value = 789.75
$ws.Range("N67").Value = -2349.75
$ws.Range("H86").Value = 1557.4166
$ws.Range("I86").Value = 1498.9
$ws.Range("K86").Value = 1498.9
$ws.Range("M86").Value = -375.9000000000001
$ws.Range("H89").Value = 1557.4166
$ws.Range("I89").Value = 1498.9
$ws.Range("K89").Value = 7494.5
$ws.Range("M89").Value = -1878.5
$ws.Range("H107").Value = 1390.5
$ws.Range("I107").Value = 1560.909
$ws.Range("K107").Value = 1560.909
$ws.Range("M107").Value = 359.0909999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1220.8
$ws.Range("I31").Value = 776
$ws.Range("K31").Value = 776
$ws.Range("M31").Value = -481
$ws.Range("H34").Value = 1220.8
$ws.Range("I34").Value = 776
$ws.Range("K34").Value = 776
$ws.Range("M34").Value = -574
$ws.Range("H58").Value = 2245.5
$ws.Range("I58").Value = 2001.375
$ws.Range("K58").Value = 2001.375
$ws.Range("M58").Value = -1798.375
$ws.Range("H86").Value = 18331.834
$ws.Range("I86").Value = 18747.75
$ws.Range("K86").Value = 18747.75
$ws.Range("M86").Value = -17624.75
$ws.Range("H89").Value = 18331.834
$ws.Range("I89").Value = 18747.75
$ws.Range("K89").Value = 93738.75
$ws.Range("M89").Value = -88122.75
$ws.Range("H99").Value = 2766.6
$ws.Range("J99").Value = 3000
$ws.Range("L99").Value = 3000
$ws.Range("N99").Value = -5996
$ws.Range("H126").Value = 2766.6
$ws.Range("J126").Value = 3000
$ws.Range("L126").Value = 9000
$ws.Range("N126").Value = -13940
$ws.Range("H136").Value = 2245.5
$ws.Range("I136").Value = 2001.375
$ws.Range("K136").Value = 6004.125
$ws.Range("M136").Value = -3454.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 176.55556
$ws.Range("I2").Value = 69.666664
$ws.Range("K2").Value = 417.999984
$ws.Range("M2").Value = -304.999984
$ws.Range("H5").Value = 333
$ws.Range("I5").Value = 333
$ws.Range("K5").Value = 999
$ws.Range("M5").Value = -887
$ws.Range("H15").Value = 806.8
$ws.Range("J15").Value = 779.6667
$ws.Range("L15").Value = 2339.0001
$ws.Range("N15").Value = -2619.0001
$ws.Range("H17").Value = 107.28571
$ws.Range("I17").Value = 91.833336
$ws.Range("K17").Value = 275.500008
$ws.Range("M17").Value = -106.500008
$ws.Range("H37").Value = 69975
$ws.Range("J37").Value = 69975
$ws.Range("L37").Value = 209925
$ws.Range("N37").Value = -210149
$ws.Range("H51").Value = 399.5
$ws.Range("I51").Value = 399.5
$ws.Range("K51").Value = 1198.5
$ws.Range("M51").Value = -738.5
$ws.Range("H106").Value = 17999.8
$ws.Range("I106").Value = 18000
$ws.Range("J106").Value = 17999
$ws.Range("K106").Value = 54000
$ws.Range("L106").Value = 53997
$ws.Range("M106").Value = -53054
$ws.Range("N106").Value = -55889
$ws.Range("H113").Value = 1698.9
$ws.Range("I113").Value = 1384.75
$ws.Range("J113").Value = 1908.3334
$ws.Range("K113").Value = 4154.25
$ws.Range("L113").Value = 5725.0002
$ws.Range("M113").Value = -1984.25
$ws.Range("N113").Value = -10065.0002
$ws.Range("H135").Value = 333
$ws.Range("I135").Value = 333
$ws.Range("K135").Value = 2997
$ws.Range("M135").Value = -462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 10624.75
$ws.Range("I10").Value = 10000
$ws.Range("J10").Value = 11249.5
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 11249.5
$ws.Range("M10").Value = -9831
$ws.Range("N10").Value = -11587.5
$ws.Range("H52").Value = 40030
$ws.Range("I52").Value = 40030
$ws.Range("K52").Value = 40030
$ws.Range("M52").Value = -39771
$ws.Range("H122").Value = 1907.625
$ws.Range("I122").Value = 1909
$ws.Range("K122").Value = 5727
$ws.Range("M122").Value = -3277
$ws.Range("H140").Value = 142773.5
$ws.Range("J140").Value = 142773.5
$ws.Range("L140").Value = 142773.5
$ws.Range("N140").Value = -153133.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2766.625
$ws.Range("I46").Value = 2304.7144
$ws.Range("J46").Value = 6000
$ws.Range("K46").Value = 2304.7144
$ws.Range("L46").Value = 6000
$ws.Range("M46").Value = -2116.7144
$ws.Range("N46").Value = -6376
$ws.Range("H74").Value = 85000
$ws.Range("I74").Value = 85000
$ws.Range("K74").Value = 85000
$ws.Range("M74").Value = -84002
$ws.Range("H77").Value = 85000
$ws.Range("I77").Value = 85000
$ws.Range("K77").Value = 255000
$ws.Range("M77").Value = -250008
$ws.Range("H125").Value = 75000
$ws.Range("J125").Value = 75000
$ws.Range("L125").Value = 75000
$ws.Range("N125").Value = -84840
$ws.Range("H132").Value = 8898
$ws.Range("I132").Value = 8898
$ws.Range("K132").Value = 26694
$ws.Range("M132").Value = -24164
$ws.Range("H136").Value = 3324.75
$ws.Range("I136").Value = 1750
$ws.Range("J136").Value = 4899.5
$ws.Range("K136").Value = 5250
$ws.Range("L136").Value = 14698.5
$ws.Range("M136").Value = -2700
$ws.Range("N136").Value = -19798.5
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2761.111
$ws.Range("I122").Value = 2761.111
$ws.Range("K122").Value = 8283.332999999999
$ws.Range("M122").Value = -5833.332999999999
